# Applies the two changes captured in the commit:
#
#  1. The table on slide 16 (the "Total Outflow / Net Cash flow" table)
#     switches its table style from
#     {349A6D5B-C2BA-4A32-BB42-47B488E3DABF} to
#     {A1D93873-1EA9-45B0-B165-6ED1F300718B}.
#
#  2. The presentation's applied theme changes from "Integral" to
#     "Office Theme" - i.e. the ten theme colours that differ between the
#     two palettes already bundled with this deck (theme1.xml / theme2.xml)
#     are updated on the active theme (theme1.xml, used by the slide
#     master/presentation) to match the "Office Theme" values.

$p = $ppt.ActivePresentation

# --- 1. Table style --------------------------------------------------------
$oldStyleId = "{349A6D5B-C2BA-4A32-BB42-47B488E3DABF}"
$newStyleId = "{A1D93873-1EA9-45B0-B165-6ED1F300718B}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable -eq -1) {
            if ($shape.Table.Style.Name -eq $oldStyleId) {
                $shape.Table.ApplyStyle($newStyleId)
            }
        }
    }
}

# --- 2. Theme colours -------------------------------------------------------
# ThemeColorScheme item order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$themeColors = $p.Slides.Item(1).ThemeColorScheme

function Set-ThemeRGB($scheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $scheme.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

Set-ThemeRGB $themeColors 1  "000000"   # dk1
Set-ThemeRGB $themeColors 2  "FFFFFF"   # lt1
Set-ThemeRGB $themeColors 3  "44546A"   # dk2
Set-ThemeRGB $themeColors 4  "E7E6E6"   # lt2
Set-ThemeRGB $themeColors 5  "5B9BD5"   # accent1
Set-ThemeRGB $themeColors 6  "ED7D31"   # accent2
Set-ThemeRGB $themeColors 7  "A5A5A5"   # accent3
Set-ThemeRGB $themeColors 8  "FFC000"   # accent4
Set-ThemeRGB $themeColors 9  "4472C4"   # accent5
Set-ThemeRGB $themeColors 10 "70AD47"   # accent6
Set-ThemeRGB $themeColors 11 "0563C1"   # hlink
Set-ThemeRGB $themeColors 12 "954F72"   # folHlink
